$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AISG Abstract")
$ws.Range("A173").Value = "test"
